$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B15").Value = "I may want to try this in Mathcad first."

$ws.Range("B15").Select()
